# Auto-generated edit script: updates pasted/refreshed market-price data
# on the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets (Kujata Profits workbook).
# Values are static data (no formulas in this workbook) so each changed
# cell is simply re-written with the new snapshot value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 738.5714
$ws.Range("I19").Value = 209.22223
$ws.Range("J19").Value = 1135.5834
$ws.Range("K19").Value = 209.22223
$ws.Range("L19").Value = 1135.5834
$ws.Range("M19").Value = -34.22223
$ws.Range("N19").Value = -1485.5834
$ws.Range("H115").Value = 4150
$ws.Range("I115").Value = 750
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 2250
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = -683
$ws.Range("N115").Value = -18134
$ws.Range("H116").Value = 2917.2104
$ws.Range("I116").Value = 2538.8147
$ws.Range("J116").Value = 3846
$ws.Range("K116").Value = 2538.8147
$ws.Range("L116").Value = 3846
$ws.Range("M116").Value = 903.1853000000001
$ws.Range("N116").Value = -10730
$ws.Range("H129").Value = 872.11
$ws.Range("J129").Value = 901.44684
$ws.Range("L129").Value = 2704.34052
$ws.Range("N129").Value = -12704.34052
$ws.Range("H132").Value = 8338025
$ws.Range("I132").Value = 9013765
$ws.Range("K132").Value = 27041295
$ws.Range("M132").Value = -27038765
$ws.Range("H138").Value = 519380.4
$ws.Range("I138").Value = 1779.2667
$ws.Range("J138").Value = 631902.4
$ws.Range("K138").Value = 5337.800099999999
$ws.Range("L138").Value = 1895707.2
$ws.Range("M138").Value = -197.8000999999995
$ws.Range("N138").Value = -1905987.2
$ws.Range("H141").Value = 325
$ws.Range("I141").Value = 325
$ws.Range("K141").Value = 975
$ws.Range("M141").Value = 4205

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1949.5
$ws.Range("I25").Value = 1949.5
$ws.Range("K25").Value = 1949.5
$ws.Range("M25").Value = -1547.5
$ws.Range("H35").Value = 2181
$ws.Range("I35").Value = 2181
$ws.Range("K35").Value = 2181
$ws.Range("M35").Value = -1775
$ws.Range("H45").Value = 1109.1538
$ws.Range("I45").Value = 1157.6666
$ws.Range("K45").Value = 1157.6666
$ws.Range("M45").Value = -780.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2506.611
$ws.Range("I31").Value = 2996.6667
$ws.Range("K31").Value = 2996.6667
$ws.Range("M31").Value = -2701.6667
$ws.Range("H34").Value = 2506.611
$ws.Range("I34").Value = 2996.6667
$ws.Range("K34").Value = 2996.6667
$ws.Range("M34").Value = -2794.6667
$ws.Range("H58").Value = 1027.08
$ws.Range("I58").Value = 733.85
$ws.Range("J58").Value = 2200
$ws.Range("K58").Value = 733.85
$ws.Range("L58").Value = 2200
$ws.Range("M58").Value = -530.85
$ws.Range("N58").Value = -2606
$ws.Range("H99").Value = 1931.25
$ws.Range("I99").Value = 1931.25
$ws.Range("K99").Value = 1931.25
$ws.Range("M99").Value = -433.25
$ws.Range("H114").Value = 23998.766
$ws.Range("J114").Value = 23998.766
$ws.Range("L114").Value = 23998.766
$ws.Range("N114").Value = -32676.766
$ws.Range("H126").Value = 1931.25
$ws.Range("I126").Value = 1931.25
$ws.Range("K126").Value = 5793.75
$ws.Range("M126").Value = -3323.75
$ws.Range("H132").Value = 7847.1055
$ws.Range("I132").Value = 10198.583
$ws.Range("K132").Value = 30595.749
$ws.Range("M132").Value = -28065.749
$ws.Range("H134").Value = 1825.0667
$ws.Range("I134").Value = 1840.3636
$ws.Range("J134").Value = 1783
$ws.Range("K134").Value = 5521.0908
$ws.Range("L134").Value = 5349
$ws.Range("M134").Value = -2986.0908
$ws.Range("N134").Value = -10419
$ws.Range("H136").Value = 1027.08
$ws.Range("I136").Value = 733.85
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 2201.55
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = 348.4499999999998
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1400.8518
$ws.Range("I5").Value = 1479.0869
$ws.Range("J5").Value = 951
$ws.Range("K5").Value = 4437.2607
$ws.Range("L5").Value = 2853
$ws.Range("M5").Value = -4325.2607
$ws.Range("N5").Value = -3077
$ws.Range("H133").Value = 4504.7144
$ws.Range("I133").Value = 1858.5
$ws.Range("J133").Value = 8033
$ws.Range("K133").Value = 5575.5
$ws.Range("L133").Value = 24099
$ws.Range("M133").Value = -515.5
$ws.Range("N133").Value = -34219
$ws.Range("H134").Value = 3765.8333
$ws.Range("I134").Value = 1845.0769
$ws.Range("J134").Value = 8759.799999999999
$ws.Range("K134").Value = 5535.2307
$ws.Range("L134").Value = 26279.4
$ws.Range("M134").Value = -465.2307000000001
$ws.Range("N134").Value = -36419.39999999999
$ws.Range("H135").Value = 1400.8518
$ws.Range("I135").Value = 1479.0869
$ws.Range("J135").Value = 951
$ws.Range("K135").Value = 13311.7821
$ws.Range("L135").Value = 8559
$ws.Range("M135").Value = -10776.7821
$ws.Range("N135").Value = -13629
$ws.Range("H137").Value = 1767.2142
$ws.Range("I137").Value = 1131.3846
$ws.Range("K137").Value = 3394.1538
$ws.Range("M137").Value = 1705.8462
$ws.Range("H139").Value = 1604.0526
$ws.Range("I139").Value = 1679.0869
$ws.Range("K139").Value = 5037.2607
$ws.Range("M139").Value = 102.7393000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6553333.5
$ws.Range("I12").Value = 6484615.5
$ws.Range("K12").Value = 6484615.5
$ws.Range("M12").Value = -6484475.5
$ws.Range("H26").Value = 13021
$ws.Range("J26").Value = 13021
$ws.Range("L26").Value = 13021
$ws.Range("N26").Value = -13581
$ws.Range("H50").Value = 13021
$ws.Range("J50").Value = 13021
$ws.Range("L50").Value = 13021
$ws.Range("N50").Value = -14017

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1880
$ws.Range("I7").Value = 1650
$ws.Range("J7").Value = 2225
$ws.Range("K7").Value = 1650
$ws.Range("L7").Value = 2225
$ws.Range("M7").Value = -1538
$ws.Range("N7").Value = -2449
$ws.Range("H122").Value = 16674953
$ws.Range("I122").Value = 22737808
$ws.Range("K122").Value = 68213424
$ws.Range("M122").Value = -68210974
$ws.Range("H126").Value = 1880
$ws.Range("I126").Value = 1650
$ws.Range("J126").Value = 2225
$ws.Range("K126").Value = 4950
$ws.Range("L126").Value = 6675
$ws.Range("M126").Value = -2480
$ws.Range("N126").Value = -11615
$ws.Range("H132").Value = 22562.562
$ws.Range("I132").Value = 1314.8334
$ws.Range("J132").Value = 57975.445
$ws.Range("K132").Value = 3944.5002
$ws.Range("L132").Value = 173926.335
$ws.Range("M132").Value = -1414.5002
$ws.Range("N132").Value = -178986.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 28667166
$ws.Range("J5").Value = 28667166
$ws.Range("L5").Value = 28667166
$ws.Range("N5").Value = -28667390
$ws.Range("H126").Value = 142858690
$ws.Range("I126").Value = 142858690
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 428576070
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -428573600
$ws.Range("N126").Value = $null
$ws.Range("H127").Value = 68200
$ws.Range("I127").Value = 55000
$ws.Range("J127").Value = 71500
$ws.Range("K127").Value = 55000
$ws.Range("L127").Value = 71500
$ws.Range("M127").Value = -50040
$ws.Range("N127").Value = -81420

